$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 155642
$ws.Range("C4").Value = 146751
$ws.Range("C5").Value = 8891
$ws.Range("C8").Value = 63.69
